$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-29"

# Update the header label cell (I1) that references the "through" date
$ws.Range("I1").Value = "2022 (through 04-29)"

# Update the May (row 5) total for the current year column (I)
$ws.Range("I5").Value = 116

# Update the grand Total row (row 14) for the current year column (I)
$ws.Range("I14").Value = 551
